$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2159.3635
$ws.Range("I80").Value = 359
$ws.Range("J80").Value = 4319.8
$ws.Range("K80").Value = 1077
$ws.Range("L80").Value = 12959.4
$ws.Range("M80").Value = -79
$ws.Range("N80").Value = -14955.4
$ws.Range("H83").Value = 2159.3635
$ws.Range("I83").Value = 359
$ws.Range("J83").Value = 4319.8
$ws.Range("K83").Value = 3231
$ws.Range("L83").Value = 38878.2
$ws.Range("M83").Value = 1761
$ws.Range("N83").Value = -48862.2
$ws.Range("H86").Value = 3090.577
$ws.Range("I86").Value = 2411.0667
$ws.Range("K86").Value = 2411.0667
$ws.Range("M86").Value = -1288.0667
$ws.Range("H88").Value = 634382.0600000001
$ws.Range("I88").Value = 1229089.9
$ws.Range("J88").Value = 89233.164
$ws.Range("K88").Value = 1229089.9
$ws.Range("L88").Value = 89233.164
$ws.Range("M88").Value = -1228683.9
$ws.Range("N88").Value = -90045.164
$ws.Range("H89").Value = 3090.577
$ws.Range("I89").Value = 2411.0667
$ws.Range("K89").Value = 12055.3335
$ws.Range("M89").Value = -6439.333499999999
$ws.Range("H91").Value = 634382.0600000001
$ws.Range("I91").Value = 1229089.9
$ws.Range("J91").Value = 89233.164
$ws.Range("K91").Value = 1229089.9
$ws.Range("L91").Value = 89233.164
$ws.Range("M91").Value = -1227685.9
$ws.Range("N91").Value = -92041.164
$ws.Range("H96").Value = 428
$ws.Range("I96").Value = 399.16666
$ws.Range("J96").Value = 514.5
$ws.Range("K96").Value = 1197.49998
$ws.Range("L96").Value = 1543.5
$ws.Range("M96").Value = 175.5000199999999
$ws.Range("N96").Value = -4289.5
$ws.Range("H97").Value = 752
$ws.Range("J97").Value = 969.3333
$ws.Range("L97").Value = 2907.9999
$ws.Range("N97").Value = -3899.9999
$ws.Range("H98").Value = 557.7241
$ws.Range("I98").Value = 557.7241
$ws.Range("K98").Value = 557.7241
$ws.Range("M98").Value = 940.2759
$ws.Range("H100").Value = 2111.7144
$ws.Range("I100").Value = 1796.4
$ws.Range("K100").Value = 1796.4
$ws.Range("M100").Value = -1255.4
$ws.Range("H111").Value = 1032.6666
$ws.Range("J111").Value = 1249.8334
$ws.Range("L111").Value = 3749.5002
$ws.Range("N111").Value = -9883.5002
$ws.Range("H122").Value = 557.7241
$ws.Range("I122").Value = 557.7241
$ws.Range("K122").Value = 1673.1723
$ws.Range("M122").Value = 776.8276999999998
$ws.Range("H132").Value = 1800.6666
$ws.Range("I132").Value = 1760.9
$ws.Range("J132").Value = 1999.5
$ws.Range("K132").Value = 5282.700000000001
$ws.Range("L132").Value = 5998.5
$ws.Range("M132").Value = -2752.700000000001
$ws.Range("N132").Value = -11058.5
$ws.Range("H137").Value = 808715.7
$ws.Range("J137").Value = 1454138.9
$ws.Range("L137").Value = 4362416.699999999
$ws.Range("N137").Value = -4367516.699999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5270.5635
$ws.Range("I32").Value = 1434.7593
$ws.Range("J32").Value = 17454.883
$ws.Range("K32").Value = 1434.7593
$ws.Range("L32").Value = 17454.883
$ws.Range("M32").Value = -1147.7593
$ws.Range("N32").Value = -18028.883
$ws.Range("H61").Value = 86223.25
$ws.Range("I61").Value = 3202.5
$ws.Range("J61").Value = 127733.625
$ws.Range("K61").Value = 3202.5
$ws.Range("L61").Value = 127733.625
$ws.Range("M61").Value = -2990.5
$ws.Range("N61").Value = -128157.625
$ws.Range("H74").Value = 34456.066
$ws.Range("I74").Value = 57083.277
$ws.Range("K74").Value = 57083.277
$ws.Range("M74").Value = -56209.277
$ws.Range("H77").Value = 34456.066
$ws.Range("I77").Value = 57083.277
$ws.Range("K77").Value = 285416.385
$ws.Range("M77").Value = -281048.385
$ws.Range("H88").Value = 1160.5294
$ws.Range("I88").Value = 1031.25
$ws.Range("J88").Value = 1275.4445
$ws.Range("K88").Value = 1031.25
$ws.Range("L88").Value = 1275.4445
$ws.Range("M88").Value = -625.25
$ws.Range("N88").Value = -2087.4445
$ws.Range("H91").Value = 1160.5294
$ws.Range("I91").Value = 1031.25
$ws.Range("J91").Value = 1275.4445
$ws.Range("K91").Value = 1031.25
$ws.Range("L91").Value = 1275.4445
$ws.Range("M91").Value = 372.75
$ws.Range("N91").Value = -4083.4445
$ws.Range("H132").Value = 2742.8572
$ws.Range("I132").Value = 2923.6667
$ws.Range("J132").Value = 2501.7778
$ws.Range("K132").Value = 8771.000100000001
$ws.Range("L132").Value = 7505.3334
$ws.Range("M132").Value = -6241.000100000001
$ws.Range("N132").Value = -12565.3334
$ws.Range("H136").Value = 86223.25
$ws.Range("I136").Value = 3202.5
$ws.Range("J136").Value = 127733.625
$ws.Range("K136").Value = 9607.5
$ws.Range("L136").Value = 383200.875
$ws.Range("M136").Value = -7057.5
$ws.Range("N136").Value = -388300.875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2369
$ws.Range("I7").Value = 216.66667
$ws.Range("J7").Value = 3983.25
$ws.Range("K7").Value = 216.66667
$ws.Range("L7").Value = 3983.25
$ws.Range("M7").Value = -103.66667
$ws.Range("N7").Value = -4209.25
$ws.Range("H99").Value = 85609.086
$ws.Range("I99").Value = 168051.5
$ws.Range("K99").Value = 168051.5
$ws.Range("M99").Value = -166553.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = $null
$ws.Range("H31").Value = 2671.3958
$ws.Range("I31").Value = 2042.0605
$ws.Range("J31").Value = 4055.9333
$ws.Range("K31").Value = 2042.0605
$ws.Range("L31").Value = 4055.9333
$ws.Range("M31").Value = -1747.0605
$ws.Range("N31").Value = -4645.933300000001
$ws.Range("H34").Value = 2671.3958
$ws.Range("I34").Value = 2042.0605
$ws.Range("J34").Value = 4055.9333
$ws.Range("K34").Value = 2042.0605
$ws.Range("L34").Value = 4055.9333
$ws.Range("M34").Value = -1840.0605
$ws.Range("N34").Value = -4459.933300000001
$ws.Range("H132").Value = 1773212.1
$ws.Range("I132").Value = 2274602.8
$ws.Range("J132").Value = 1001841.75
$ws.Range("K132").Value = 6823808.399999999
$ws.Range("L132").Value = 3005525.25
$ws.Range("M132").Value = -6821278.399999999
$ws.Range("N132").Value = -3010585.25
$ws.Range("H134").Value = 3020692.8
$ws.Range("I134").Value = 5956023.5
$ws.Range("J134").Value = 85361.914
$ws.Range("K134").Value = 17868070.5
$ws.Range("L134").Value = 256085.742
$ws.Range("M134").Value = -17865535.5
$ws.Range("N134").Value = -261155.742

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 6723.5835
$ws.Range("I132").Value = 1995
$ws.Range("J132").Value = 7153.4546
$ws.Range("K132").Value = 17955
$ws.Range("L132").Value = 64381.0914
$ws.Range("M132").Value = -15425
$ws.Range("N132").Value = -69441.0914

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 499.2
$ws.Range("I5").Value = 499.2
$ws.Range("K5").Value = 499.2
$ws.Range("M5").Value = -387.2
$ws.Range("H80").Value = 3149.25
$ws.Range("I80").Value = 2599.4
$ws.Range("K80").Value = 2599.4
$ws.Range("M80").Value = -1601.4
$ws.Range("H83").Value = 3149.25
$ws.Range("I83").Value = 2599.4
$ws.Range("K83").Value = 12997
$ws.Range("M83").Value = -8005
$ws.Range("H113").Value = 3589.125
$ws.Range("I113").Value = 3237
$ws.Range("J113").Value = 3800.4
$ws.Range("K113").Value = 3237
$ws.Range("L113").Value = 3800.4
$ws.Range("M113").Value = -1067
$ws.Range("N113").Value = -8140.4
$ws.Range("H120").Value = 44500
$ws.Range("J120").Value = 44500
$ws.Range("L120").Value = 44500
$ws.Range("N120").Value = -54176
$ws.Range("H132").Value = 7104.6
$ws.Range("I132").Value = 8130.75
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 24392.25
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -21862.25
$ws.Range("N132").Value = -14060

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1169.7142
$ws.Range("I61").Value = 1169.7142
$ws.Range("K61").Value = 1169.7142
$ws.Range("M61").Value = -967.7141999999999
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null
$ws.Range("H100").Value = 4893.4736
$ws.Range("I100").Value = 4292.7646
$ws.Range("K100").Value = 4292.7646
$ws.Range("M100").Value = -3751.7646
$ws.Range("H113").Value = 1169.7142
$ws.Range("I113").Value = 1169.7142
$ws.Range("K113").Value = 1169.7142
$ws.Range("M113").Value = 1000.2858
$ws.Range("H136").Value = 2098.9167
$ws.Range("I136").Value = 1585.5
$ws.Range("J136").Value = 2817.7
$ws.Range("K136").Value = 4756.5
$ws.Range("L136").Value = 8453.099999999999
$ws.Range("M136").Value = -2206.5
$ws.Range("N136").Value = -13553.1

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 34833
$ws.Range("J70").Value = 37799.6
$ws.Range("L70").Value = 37799.6
$ws.Range("N70").Value = -38429.6
$ws.Range("H73").Value = 34833
$ws.Range("J73").Value = 37799.6
$ws.Range("L73").Value = 37799.6
$ws.Range("N73").Value = -39983.6
$ws.Range("H75").Value = 41600
$ws.Range("I75").Value = 23333.334
$ws.Range("K75").Value = 23333.334
$ws.Range("M75").Value = -22397.334
$ws.Range("H78").Value = 41600
$ws.Range("I78").Value = 23333.334
$ws.Range("K78").Value = 70000.00199999999
$ws.Range("M78").Value = -65320.00199999999
$ws.Range("H80").Value = 72500
$ws.Range("J80").Value = 72500
$ws.Range("L80").Value = 72500
$ws.Range("N80").Value = -74496
$ws.Range("H83").Value = 72500
$ws.Range("J83").Value = 72500
$ws.Range("L83").Value = 217500
$ws.Range("N83").Value = -227484
$ws.Range("H132").Value = 2092.4
$ws.Range("I132").Value = 1810.8422
$ws.Range("K132").Value = 5432.5266
$ws.Range("M132").Value = -2902.5266
